# Apply the "solvents" sheet update: add a new "vapour pressure" column (F)
# with per-solvent values, apply wrap-text formatting to one cell, and move
# the active-sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solvents")

# Values entered in the same order the author typed them (this determines
# the order new entries land in the shared-strings table)
$ws.Range("F3").Value  = "516 Pa"       # Dimethyl formamide
$ws.Range("F1").Value  = "vapour pressure"   # header
$ws.Range("F2").Value  = "200 mmHg"     # Chloroform
$ws.Range("F4").Value  = "115 mmHg"     # Hexafluoro isopropanol
$ws.Range("F6").Value  = "24 kPa"       # Acetone
$ws.Range("F7").Value  = "760 mmHg"     # Water
$ws.Range("F8").Value  = "13.02 kPa"    # Methanol
$ws.Range("F9").Value  = "202 kPa"      # Acetic acid
$ws.Range("F10").Value = "4.6 kPa"      # Formic acid
$ws.Range("F11").Value = "53.3 kPa"     # Dichloro methane
$ws.Range("F12").Value = "5.95 kPa"     # Ethanol
$ws.Range("F13").Value = "96.2 mmHg"    # Tri fluoro acetic acid
$ws.Range("F14").Value = "0.42 mmHg"    # Dimethyl sulfoxide

# Hexafluoro isopropanol's vapour pressure cell is wrapped
$ws.Range("F4").WrapText = $true

# Give the new column a sensible best-fit width
$ws.Columns.Item(6).ColumnWidth = 14.6

# Move the selection on the sheet (matches the saved cursor position)
$ws.Range("F25").Select() | Out-Null
